$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 4 (old row4 -> row7, old row12 -> row15;
# rows 1-3 stay put for now).
$ws.Rows("4:6").Insert()

# Row 3 (original "preproduccion..." claim, with the hyperlink) still needs
# to be duplicated into the new rows 5 and 6 before it gets overwritten.
# (Copy bounded A:G ranges, not whole rows, to avoid materializing every
# column out to XFD.)
$ws.Range("A3:G3").Copy($ws.Range("A5:G5"))
$ws.Range("A3:G3").Copy($ws.Range("A6:G6"))

# Row 2 (original "ssurgwsoadev..." claim) gets duplicated into rows 3 and 4.
$ws.Range("A2:G2").Copy($ws.Range("A3:G3"))
$ws.Range("A2:G2").Copy($ws.Range("A4:G4"))

# The single hyperlink that used to live on B3 needs to be re-anchored: one
# copy stays on the row that keeps the original relationship (B6) and a new
# hyperlink relationship is added for the duplicate on B5.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
# Hyperlinks.Add() re-derives cell formatting; snap both cells back onto the
# existing "Hipervinculo" cell style (they already carry it from the row
# copy above) instead of leaving a duplicate equivalent style behind.
$ws.Range("B5").Style = "Hipervínculo"
$ws.Range("B6").Style = "Hipervínculo"

# Row 2 itself changes to a new claim (quita acentos de los TCs): new user
# and a brand-new claim number. The NroSiniestro column is formatted as
# Text with a quote-prefix (so Excel doesn't try to treat the long numeric
# string as a number) - enter it with a leading apostrophe so that marker
# is preserved, same as the original cell.
$ws.Range("C2").Value = "apellegrini"
$ws.Range("F2").Value = "'1220194200662"

# Update the saved selection.
$ws.Range("C3").Select()
